$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 641440.1
$ws.Range("I64").Value = 1460809.2
$ws.Range("J64").Value = 4153
$ws.Range("K64").Value = 1460809.2
$ws.Range("L64").Value = 4153
$ws.Range("M64").Value = -1460561.2
$ws.Range("N64").Value = -4649
$ws.Range("H67").Value = 641440.1
$ws.Range("I67").Value = 1460809.2
$ws.Range("J67").Value = 4153
$ws.Range("K67").Value = 1460809.2
$ws.Range("L67").Value = 4153
$ws.Range("M67").Value = -1459951.2
$ws.Range("N67").Value = -5869
$ws.Range("H129").Value = 1082.5634
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670
$ws.Range("H32").Value = 19638.42
$ws.Range("I32").Value = 22029.76
$ws.Range("K32").Value = 22029.76
$ws.Range("M32").Value = -21742.76
$ws.Range("H122").Value = 1949.2142
$ws.Range("I122").Value = 1717.1818
$ws.Range("K122").Value = 5151.5454
$ws.Range("M122").Value = -2701.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 50850
$ws.Range("I59").Value = 35000
$ws.Range("J59").Value = 58775
$ws.Range("K59").Value = 35000
$ws.Range("L59").Value = 58775
$ws.Range("M59").Value = -34153
$ws.Range("N59").Value = -60469
$ws.Range("H105").Value = 2211.111
$ws.Range("I105").Value = 2616.6667
$ws.Range("J105").Value = 1400
$ws.Range("K105").Value = 2616.6667
$ws.Range("L105").Value = 1400
$ws.Range("M105").Value = -869.6667000000002
$ws.Range("N105").Value = -4894
$ws.Range("H134").Value = 2337.5483
$ws.Range("I134").Value = 2215.2
$ws.Range("J134").Value = 2847.3333
$ws.Range("K134").Value = 6645.599999999999
$ws.Range("L134").Value = 8541.999899999999
$ws.Range("M134").Value = -4110.599999999999
$ws.Range("N134").Value = -13611.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 241.07692
$ws.Range("I22").Value = 221
$ws.Range("J22").Value = 264.5
$ws.Range("K22").Value = 221
$ws.Range("L22").Value = 264.5
$ws.Range("M22").Value = 129
$ws.Range("N22").Value = -964.5
$ws.Range("H62").Value = 3149.5
$ws.Range("I62").Value = 2979.4
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 2979.4
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2355.4
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3149.5
$ws.Range("I65").Value = 2979.4
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 14897
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -11777
$ws.Range("N65").Value = -26240
$ws.Range("H80").Value = 16000
$ws.Range("J80").Value = 16000
$ws.Range("L80").Value = 16000
$ws.Range("N80").Value = -18246
$ws.Range("H83").Value = 16000
$ws.Range("J83").Value = 16000
$ws.Range("L83").Value = 48000
$ws.Range("N83").Value = -59232
$ws.Range("H122").Value = 7391.304
$ws.Range("I122").Value = 4316.857
$ws.Range("J122").Value = 12173.777
$ws.Range("K122").Value = 12950.571
$ws.Range("L122").Value = 36521.331
$ws.Range("M122").Value = -10500.571
$ws.Range("N122").Value = -41421.331
$ws.Range("H132").Value = 3425.3635
$ws.Range("I132").Value = 3537
$ws.Range("J132").Value = 3098.4285
$ws.Range("K132").Value = 10611
$ws.Range("L132").Value = 9295.2855
$ws.Range("M132").Value = -8081
$ws.Range("N132").Value = -14355.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 493.375
$ws.Range("I98").Value = 391.5
$ws.Range("K98").Value = 1174.5
$ws.Range("M98").Value = 323.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3338.6
$ws.Range("I122").Value = 5004.769
$ws.Range("J122").Value = 1533.5834
$ws.Range("K122").Value = 15014.307
$ws.Range("L122").Value = 4600.7502
$ws.Range("M122").Value = -12564.307
$ws.Range("N122").Value = -9500.7502
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -69820
$ws.Range("H132").Value = 2266.875
$ws.Range("I132").Value = 2158.348
$ws.Range("J132").Value = 2544.2222
$ws.Range("K132").Value = 6475.044
$ws.Range("L132").Value = 7632.6666
$ws.Range("M132").Value = -3945.044
$ws.Range("N132").Value = -12692.6666
$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -45120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5680
$ws.Range("I7").Value = 4906.6665
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 4906.6665
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -4794.6665
$ws.Range("N7").Value = -8224
$ws.Range("H22").Value = 957
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 935.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 935.5
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1525.5
$ws.Range("H27").Value = 957
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 935.5
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 935.5
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1149.5
$ws.Range("H93").Value = 1000.58826
$ws.Range("J93").Value = 1270.5
$ws.Range("L93").Value = 1270.5
$ws.Range("N93").Value = -3766.5
$ws.Range("H126").Value = 5680
$ws.Range("I126").Value = 4906.6665
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 14719.9995
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -12249.9995
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 2531.8438
$ws.Range("I132").Value = 1964.84
$ws.Range("K132").Value = 5894.52
$ws.Range("M132").Value = -3364.52

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1799.907
$ws.Range("I132").Value = 978.62964
$ws.Range("K132").Value = 2935.88892
$ws.Range("M132").Value = -405.8889199999999
